$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9033119704407397
$ws.Range("C2").Value = 0.04477016915814147
$ws.Range("D2").Value = 0.1068254417072581
$ws.Range("F2").Value = 2.1521882487867
$ws.Range("G2").Value = 1.476200289105861
$ws.Range("H2").Value = 1.363154875329542
$ws.Range("I2").Value = 1.343990568854402
$ws.Range("J2").Value = 0.1849649801769573
$ws.Range("K2").Value = 0.6014323791593483
$ws.Range("L2").Value = 0.3808079546235632
$ws.Range("N2").Value = 2.348950722251576
$ws.Range("B3").Value = 0.8655104710579167
$ws.Range("C3").Value = 0.04141837103514945
$ws.Range("D3").Value = 0.1053583587281821
$ws.Range("F3").Value = 2.156619209380885
$ws.Range("G3").Value = 1.479313592385054
$ws.Range("H3").Value = 1.36933237111856
$ws.Range("I3").Value = 1.350655206317597
$ws.Range("J3").Value = 0.1857199811169394
$ws.Range("K3").Value = 0.564087855340631
$ws.Range("L3").Value = 0.3761041292347471
$ws.Range("N3").Value = 2.370117659928827
$ws.Range("B4").Value = 0.8426876459120081
$ws.Range("C4").Value = 0.03934054963306011
$ws.Range("D4").Value = 0.1044929233538952
$ws.Range("F4").Value = 2.160275204189951
$ws.Range("G4").Value = 1.481949645443265
$ws.Range("H4").Value = 1.373627914101277
$ws.Range("I4").Value = 1.35529283058959
$ws.Range("J4").Value = 0.1862398218708101
$ws.Range("K4").Value = 0.5414019027421659
$ws.Range("L4").Value = 0.3733773079998599
$ws.Range("N4").Value = 2.383786024350385
$ws.Range("B5").Value = 0.8334852210069243
$ws.Range("C5").Value = 0.03848883345310128
$ws.Range("D5").Value = 0.1041491971535109
$ws.Range("F5").Value = 2.162000469242713
$ws.Range("G5").Value = 1.483206078335357
$ws.Range("H5").Value = 1.375504879375725
$ws.Range("I5").Value = 1.357319962875142
$ws.Range("J5").Value = 0.1864658353817887
$ws.Range("K5").Value = 0.532218954883092
$ws.Range("L5").Value = 0.3723068097506967
$ws.Range("N5").Value = 2.389524940731482
$ws.Range("B6").Value = 0.831963106746997
$ws.Range("C6").Value = 0.03834710521317675
$ws.Range("D6").Value = 0.1040926635790527
$ws.Range("F6").Value = 2.162301173086
$ws.Range("G6").Value = 1.483425715057791
$ws.Range("H6").Value = 1.37582419195283
$ws.Range("I6").Value = 1.357664859637659
$ws.Range("J6").Value = 0.1865042215858423
$ws.Range("K6").Value = 0.5306978775978735
$ws.Range("L6").Value = 0.3721315166039716
$ws.Range("N6").Value = 2.390488088110178
$ws.Range("B7").Value = 0.8425631409026835
$ws.Range("C7").Value = 0.0393290832931541
$ws.Range("D7").Value = 0.1044882514489274
$ws.Range("F7").Value = 2.160297518242068
$ws.Range("G7").Value = 1.481965852276375
$ws.Range("H7").Value = 1.373652715161242
$ws.Range("I7").Value = 1.355319613315146
$ws.Range("J7").Value = 0.1862428125432984
$ws.Range("K7").Value = 0.541277807596714
$ws.Range("L7").Value = 0.3733627059040572
$ws.Range("N7").Value = 2.383862737296434
$ws.Range("B8").Value = 0.8901979655090884
$ws.Range("C8").Value = 0.04361857420049375
$ws.Range("D8").Value = 0.1063122845329403
$ws.Range("F8").Value = 2.153522031751407
$ws.Range("G8").Value = 1.477123416073397
$ws.Range("H8").Value = 1.365180651504801
$ws.Range("I8").Value = 1.346175368060067
$ws.Range("J8").Value = 0.1852136413684811
$ws.Range("K8").Value = 0.5885056707578258
$ws.Range("L8").Value = 0.3791526647936507
$ws.Range("N8").Value = 2.356109670093193
$ws.Range("B9").Value = 0.9866615746891227
$ws.Range("C9").Value = 0.05187380685931942
$ws.Range("D9").Value = 0.1101676436615548
$ws.Range("F9").Value = 2.147649479304732
$ws.Range("G9").Value = 1.473375201917491
$ws.Range("H9").Value = 1.35254912558635
$ws.Range("I9").Value = 1.332568490973287
$ws.Range("J9").Value = 0.1836408223232233
$ws.Range("K9").Value = 0.68303767119545
$ws.Range("L9").Value = 0.3917821372305639
$ws.Range("N9").Value = 2.307014434897649
$ws.Range("B10").Value = 1.05937196354904
$ws.Range("C10").Value = 0.05784502467405161
$ws.Range("D10").Value = 0.1131674427970495
$ws.Range("F10").Value = 2.147847201670288
$ws.Range("G10").Value = 1.474126734856284
$ws.Range("H10").Value = 1.345690061234535
$ws.Range("I10").Value = 1.325204325092926
$ws.Range("J10").Value = 0.182755410505429
$ws.Range("K10").Value = 0.7536463492499763
$ws.Range("L10").Value = 0.4018334262373884
$ws.Range("N10").Value = 2.27418886012422
$ws.Range("B11").Value = 1.092844885654642
$ws.Range("C11").Value = 0.06054144074717271
$ws.Range("D11").Value = 0.1145679495678849
$ws.Range("F11").Value = 2.148915546162044
$ws.Range("G11").Value = 1.475230229271531
$ws.Range("H11").Value = 1.343094196750371
$ws.Range("I11").Value = 1.322425136950748
$ws.Range("J11").Value = 0.1824109861601393
$ws.Range("K11").Value = 0.7860166966557074
$ws.Range("L11").Value = 0.4065727507695414
$ws.Range("N11").Value = 2.259959434796869
$ws.Range("B12").Value = 1.105576711917934
$ws.Range("C12").Value = 0.06155965395076635
$ws.Range("D12").Value = 0.115103395295364
$ws.Range("F12").Value = 2.149460634805635
$ws.Range("G12").Value = 1.475757616162284
$ws.Range("H12").Value = 1.342186502456457
$ws.Range("I12").Value = 1.321454745327586
$ws.Range("J12").Value = 0.1822889289368526
$ws.Range("K12").Value = 0.7983100965282972
$ws.Range("L12").Value = 0.40839130337028
$ws.Range("N12").Value = 2.254672211811016
$ws.Range("B13").Value = 1.1028321895555
$ws.Range("C13").Value = 0.06134049070993797
$ws.Range("D13").Value = 0.1149878513524669
$ws.Range("F13").Value = 2.149336992699233
$ws.Range("G13").Value = 1.475639163170996
$ws.Range("H13").Value = 1.342378643093653
$ws.Range("I13").Value = 1.321660089236069
$ws.Range("J13").Value = 0.1823148442946056
$ws.Range("K13").Value = 0.7956609243491073
$ws.Range("L13").Value = 0.4079985858758874
$ws.Range("N13").Value = 2.255806410345947
$ws.Range("B14").Value = 1.093891215041111
$ws.Range("C14").Value = 0.06062526708002736
$ws.Range("D14").Value = 0.1146118990329512
$ws.Range("F14").Value = 2.148957575437237
$ws.Range("G14").Value = 1.475271422826452
$ws.Range("H14").Value = 1.343018011567892
$ws.Range("I14").Value = 1.322343658624121
$ws.Range("J14").Value = 0.1824007768434441
$ws.Range("K14").Value = 0.7870273744940448
$ws.Range("L14").Value = 0.40672188649188
$ws.Range("N14").Value = 2.259522424268088
$ws.Range("B15").Value = 1.088421929357281
$ws.Range("C15").Value = 0.06018679946673444
$ws.Range("D15").Value = 0.1143822804485808
$ws.Range("F15").Value = 2.148743467711668
$ws.Range("G15").Value = 1.475060433597136
$ws.Range("H15").Value = 1.343419447008813
$ws.Range("I15").Value = 1.322773045282084
$ws.Range("J15").Value = 0.1824545022121775
$ws.Range("K15").Value = 0.7817436775646911
$ws.Range("L15").Value = 0.4059429763150035
$ws.Range("N15").Value = 2.261811764479607
$ws.Range("B16").Value = 1.057192352407839
$ws.Range("C16").Value = 0.05766840698929343
$ws.Range("D16").Value = 0.1130766336780624
$ws.Range("F16").Value = 2.147797056803142
$ws.Range("G16").Value = 1.474069944405258
$ws.Range("H16").Value = 1.345870250041614
$ws.Range("I16").Value = 1.325397432980239
$ws.Range("J16").Value = 0.1827790917570944
$ws.Range("K16").Value = 0.7515358598305966
$ws.Range("L16").Value = 0.4015270481853861
$ws.Range("N16").Value = 2.275132938564862
$ws.Range("B17").Value = 1.038135118191292
$ws.Range("C17").Value = 0.05611835513019514
$ws.Range("D17").Value = 0.1122848124968385
$ws.Range("F17").Value = 2.147466939820688
$ws.Range("G17").Value = 1.473657378358837
$ws.Range("H17").Value = 1.347507968425219
$ws.Range("I17").Value = 1.327153568143039
$ws.Range("J17").Value = 0.1829931467795873
$ws.Range("K17").Value = 0.7330680351450951
$ws.Range("L17").Value = 0.3988606827630576
$ws.Range("N17").Value = 2.283485159320012
$ws.Range("B18").Value = 1.027211271377524
$ws.Range("C18").Value = 0.05522493104628268
$ws.Range("D18").Value = 0.1118327587095251
$ws.Range("F18").Value = 2.147369166543115
$ws.Range("G18").Value = 1.4734917685075
$ws.Range("H18").Value = 1.348499301644438
$ws.Range("I18").Value = 1.32821738100246
$ws.Range("J18").Value = 0.1831217594482908
$ws.Range("K18").Value = 0.722469402644947
$ws.Range("L18").Value = 0.3973427812789367
$ws.Range("N18").Value = 2.288355311194493
$ws.Range("B19").Value = 1.023519083455909
$ws.Range("C19").Value = 0.05492211081771359
$ws.Range("D19").Value = 0.1116802834213573
$ws.Range("F19").Value = 2.147351885655738
$ws.Range("G19").Value = 1.473448008286226
$ws.Range("H19").Value = 1.348843430855794
$ws.Range("I19").Value = 1.328586799686349
$ws.Range("J19").Value = 0.1831662497527127
$ws.Range("K19").Value = 0.7188849518999803
$ws.Range("L19").Value = 0.3968315501312105
$ws.Range("N19").Value = 2.290015625234926
$ws.Range("B20").Value = 1.040159930505666
$ws.Range("C20").Value = 0.05628355485235659
$ws.Range("D20").Value = 0.1123687537172131
$ws.Range("F20").Value = 2.147492550380079
$ws.Range("G20").Value = 1.473693877375212
$ws.Range("H20").Value = 1.347328522553781
$ws.Range("I20").Value = 1.326961064063838
$ws.Range("J20").Value = 0.1829697918237585
$ws.Range("K20").Value = 0.7350315307329254
$ws.Range("L20").Value = 0.3991428957737071
$ws.Range("N20").Value = 2.282589201917208
$ws.Range("B21").Value = 1.096515871840495
$ws.Range("C21").Value = 0.06083542321761115
$ws.Range("D21").Value = 0.1147221872588915
$ws.Range("F21").Value = 2.149065206840945
$ws.Range("G21").Value = 1.475376464822745
$ws.Range("H21").Value = 1.342828170622511
$ws.Range("I21").Value = 1.322140651979652
$ws.Range("J21").Value = 0.1823753094220031
$ws.Range("K21").Value = 0.7895622993034976
$ws.Range("L21").Value = 0.4070962371943097
$ws.Range("N21").Value = 2.258428195281898
$ws.Range("B22").Value = 1.133675917322932
$ws.Range("C22").Value = 0.06379365754661137
$ws.Range("D22").Value = 0.1162900175810648
$ws.Range("F22").Value = 2.150912068764313
$ws.Range("G22").Value = 1.477114484003621
$ws.Range("H22").Value = 1.340325810277534
$ws.Range("I22").Value = 1.319468322951913
$ws.Range("J22").Value = 0.1820355502646827
$ws.Range("K22").Value = 0.8254077111070046
$ws.Range("L22").Value = 0.4124332955654211
$ws.Range("N22").Value = 2.243227193827774
$ws.Range("B23").Value = 1.113813095721525
$ws.Range("C23").Value = 0.06221631680973871
$ws.Range("D23").Value = 0.1154505355336113
$ws.Range("F23").Value = 2.149851475204969
$ws.Range("G23").Value = 1.476128461963484
$ws.Range("H23").Value = 1.341621241960581
$ws.Range("I23").Value = 1.320850868015555
$ws.Range("J23").Value = 0.1822124310904556
$ws.Range("K23").Value = 0.8062576239622672
$ws.Range("L23").Value = 0.4095721242197641
$ws.Range("N23").Value = 2.251286288203833
$ws.Range("B24").Value = 1.039244411850348
$ws.Range("C24").Value = 0.05620887515546258
$ws.Range("D24").Value = 0.1123307939985807
$ws.Range("F24").Value = 2.147480685208919
$ws.Range("G24").Value = 1.473677153200285
$ws.Range("H24").Value = 1.347409494937509
$ws.Range("I24").Value = 1.327047926360862
$ws.Range("J24").Value = 0.1829803333110682
$ws.Range("K24").Value = 0.734143775885741
$ws.Range("L24").Value = 0.3990152604536235
$ws.Range("N24").Value = 2.282994051352915
$ws.Range("B25").Value = 0.9602411579113834
$ws.Range("C25").Value = 0.04965714377040342
$ws.Range("D25").Value = 0.1090951147737371
$ws.Range("F25").Value = 2.148445211611801
$ws.Range("G25").Value = 1.473773716584518
$ws.Range("H25").Value = 1.355540626963176
$ws.Range("I25").Value = 1.33578689167274
$ws.Range("J25").Value = 0.1840187732844321
$ws.Range("K25").Value = 0.6572602814439392
$ws.Range("L25").Value = 0.3882295187404026
$ws.Range("N25").Value = 2.31972594193947
